$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 1; $r -le 10; $r++) {
    # Copy the style/formatting from the adjacent phone-number cell (column F)
    # onto the new account-number cell (column G) so it picks up the same
    # cellXf (s="3") used by the other data columns.
    $ws.Cells.Item($r, 6).Copy($ws.Cells.Item($r, 7))

    $phone = $ws.Cells.Item($r, 6).Value2
    $tail = $phone.Substring(2)
    $accountNumber = "0$tail"

    $ws.Cells.Item($r, 7).Value = $accountNumber
}
